$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.514.71"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.449.51"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.91"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.09"
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("D7").Value = "3.441.29"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.06"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.565"
$ws.Range("E12").Value = "  -2.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.39"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000270"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "4.003.01"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "583.37"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("D18").Value = "69.625.08"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "3.447.32"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.91"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.00"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.21"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.85"
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.66"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.56"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "574.08"
$ws.Range("E35").Value = "  -13.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0480"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.14"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("E42").Value = "  -10.68%  "
$ws.Range("D43").Value = "3.244.94"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "0.0₃0687"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.22"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.296"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  -4.98%  "
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.40"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +0.01%  "
